$d = $word.ActiveDocument

# First paragraph of the document holds the placeholder ID text.
$para = $d.Paragraphs(1)

# Add a thin paragraph border (5-twip space) on all four sides.
$borders = $para.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Increase the left indent of the paragraph from 120 to 225 twips
# (ParagraphFormat.LeftIndent is expressed in points => 225/20 = 11.25pt).
$para.Format.LeftIndent = 11.25

# Replace the placeholder text (together with its trailing-space run) with
# the new ID, collapsing the paragraph back down to a single run.
# MatchWildcards is left $false so the literal asterisks in the ID text are
# not interpreted as wildcard metacharacters.
$d.Content.Find.Execute("**ID__AFFARS_5313_topic_6__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5313_201__ID**", 2)
